# Update TPM-derived NATMI ligand-receptor metrics (Sema4a-Plxnb1) with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 8.050906
$ws.Range("H2").Value2 = 24.152718
$ws.Range("I2").Value2 = 0.1888708516018927
$ws.Range("J2").Value2 = 0.1888708516018927
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 1.480335666666667
$ws.Range("N2").Value2 = 4.441007
$ws.Range("O2").Value2 = 0.1826408776454046
$ws.Range("P2").Value2 = 0.1826408776454046
$ws.Range("Q2").Value2 = 11.91804330078067
$ws.Range("R2").Value2 = 107.262389707026
$ws.Range("S2").Value2 = 0.03449553809820465
$ws.Range("T2").Value2 = 0.03449553809820464

$ws.Range("G3").Value2 = 8.050906
$ws.Range("H3").Value2 = 24.152718
$ws.Range("I3").Value2 = 0.1888708516018927
$ws.Range("J3").Value2 = 0.1888708516018927
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 0.6291593333333333
$ws.Range("N3").Value2 = 1.887478
$ws.Range("O3").Value2 = 0.07762443032771463
$ws.Range("P3").Value2 = 0.07762443032771463
$ws.Range("Q3").Value2 = 5.065302651689333
$ws.Range("R3").Value2 = 45.587723865204
$ws.Range("S3").Value2 = 0.01466099226110724
$ws.Range("T3").Value2 = 0.01466099226110724

$ws.Range("G4").Value2 = 8.050906
$ws.Range("H4").Value2 = 24.152718
$ws.Range("I4").Value2 = 0.1888708516018927
$ws.Range("J4").Value2 = 0.1888708516018927
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 5.953764333333333
$ws.Range("N4").Value2 = 17.861293
$ws.Range("O4").Value2 = 0.734563631492074
$ws.Range("P4").Value2 = 0.734563631492074
$ws.Range("Q4").Value2 = 47.93319699381933
$ws.Range("R4").Value2 = 431.398772944374
$ws.Range("S4").Value2 = 0.1387376586356869
$ws.Range("T4").Value2 = 0.1387376586356869

$ws.Range("G5").Value2 = 8.050906
$ws.Range("H5").Value2 = 24.152718
$ws.Range("I5").Value2 = 0.1888708516018927
$ws.Range("J5").Value2 = 0.1888708516018927
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.04191233333333333
$ws.Range("N5").Value2 = 0.125737
$ws.Range("O5").Value2 = 0.005171060534806686
$ws.Range("P5").Value2 = 0.005171060534806686
$ws.Range("Q5").Value2 = 0.3374322559073333
$ws.Range("R5").Value2 = 3.036890303166
$ws.Range("S5").Value2 = 0.0009766626068938773
$ws.Range("T5").Value2 = 0.0009766626068938773

$ws.Range("G6").Value2 = 12.38193366666667
$ws.Range("H6").Value2 = 37.145801
$ws.Range("I6").Value2 = 0.2904749299149038
$ws.Range("J6").Value2 = 0.2904749299149038
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 1.480335666666667
$ws.Range("N6").Value2 = 4.441007
$ws.Range("O6").Value2 = 0.1826408776454046
$ws.Range("P6").Value2 = 0.1826408776454046
$ws.Range("Q6").Value2 = 18.32941802906744
$ws.Range("R6").Value2 = 164.964762261607
$ws.Range("S6").Value2 = 0.05305259613364542
$ws.Range("T6").Value2 = 0.05305259613364541

$ws.Range("G7").Value2 = 12.38193366666667
$ws.Range("H7").Value2 = 37.145801
$ws.Range("I7").Value2 = 0.2904749299149038
$ws.Range("J7").Value2 = 0.2904749299149038
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 0.6291593333333333
$ws.Range("N7").Value2 = 1.887478
$ws.Range("O7").Value2 = 0.07762443032771463
$ws.Range("P7").Value2 = 0.07762443032771463
$ws.Range("Q7").Value2 = 7.790209131097555
$ws.Range("R7").Value2 = 70.111882179878
$ws.Range("S7").Value2 = 0.02254795095912724
$ws.Range("T7").Value2 = 0.02254795095912724

$ws.Range("G8").Value2 = 12.38193366666667
$ws.Range("H8").Value2 = 37.145801
$ws.Range("I8").Value2 = 0.2904749299149038
$ws.Range("J8").Value2 = 0.2904749299149038
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 5.953764333333333
$ws.Range("N8").Value2 = 17.861293
$ws.Range("O8").Value2 = 0.734563631492074
$ws.Range("P8").Value2 = 0.734563631492074
$ws.Range("Q8").Value2 = 73.71911504229922
$ws.Range("R8").Value2 = 663.472035380693
$ws.Range("S8").Value2 = 0.2133723193756975
$ws.Range("T8").Value2 = 0.2133723193756975

$ws.Range("G9").Value2 = 12.38193366666667
$ws.Range("H9").Value2 = 37.145801
$ws.Range("I9").Value2 = 0.2904749299149038
$ws.Range("J9").Value2 = 0.2904749299149038
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.04191233333333333
$ws.Range("N9").Value2 = 0.125737
$ws.Range("O9").Value2 = 0.005171060534806686
$ws.Range("P9").Value2 = 0.005171060534806686
$ws.Range("Q9").Value2 = 0.5189557311485555
$ws.Range("R9").Value2 = 4.670601580336999
$ws.Range("S9").Value2 = 0.001502063446433697
$ws.Range("T9").Value2 = 0.001502063446433697

$ws.Range("G10").Value2 = 8.230170000000001
$ws.Range("H10").Value2 = 24.69051
$ws.Range("I10").Value2 = 0.1930763092661061
$ws.Range("J10").Value2 = 0.1930763092661061
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 1.480335666666667
$ws.Range("N10").Value2 = 4.441007
$ws.Range("O10").Value2 = 0.1826408776454046
$ws.Range("P10").Value2 = 0.1826408776454046
$ws.Range("Q10").Value2 = 12.18341419373
$ws.Range("R10").Value2 = 109.65072774357
$ws.Range("S10").Value2 = 0.03526362657689718
$ws.Range("T10").Value2 = 0.03526362657689718

$ws.Range("G11").Value2 = 8.230170000000001
$ws.Range("H11").Value2 = 24.69051
$ws.Range("I11").Value2 = 0.1930763092661061
$ws.Range("J11").Value2 = 0.1930763092661061
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 0.6291593333333333
$ws.Range("N11").Value2 = 1.887478
$ws.Range("O11").Value2 = 0.07762443032771463
$ws.Range("P11").Value2 = 0.07762443032771463
$ws.Range("Q11").Value2 = 5.17808827042
$ws.Range("R11").Value2 = 46.60279443378001
$ws.Range("S11").Value2 = 0.01498743851655913
$ws.Range("T11").Value2 = 0.01498743851655913

$ws.Range("G12").Value2 = 8.230170000000001
$ws.Range("H12").Value2 = 24.69051
$ws.Range("I12").Value2 = 0.1930763092661061
$ws.Range("J12").Value2 = 0.1930763092661061
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 5.953764333333333
$ws.Range("N12").Value2 = 17.861293
$ws.Range("O12").Value2 = 0.734563631492074
$ws.Range("P12").Value2 = 0.734563631492074
$ws.Range("Q12").Value2 = 49.00049260327
$ws.Range("R12").Value2 = 441.00443342943
$ws.Range("S12").Value2 = 0.1418268348895977
$ws.Range("T12").Value2 = 0.1418268348895977

$ws.Range("G13").Value2 = 8.230170000000001
$ws.Range("H13").Value2 = 24.69051
$ws.Range("I13").Value2 = 0.1930763092661061
$ws.Range("J13").Value2 = 0.1930763092661061
$ws.Range("K13").Value2 = 1
$ws.Range("L13").Value2 = 0.3333333333333333
$ws.Range("M13").Value2 = 0.04191233333333333
$ws.Range("N13").Value2 = 0.125737
$ws.Range("O13").Value2 = 0.005171060534806686
$ws.Range("P13").Value2 = 0.005171060534806686
$ws.Range("Q13").Value2 = 0.34494562843
$ws.Range("R13").Value2 = 3.10451065587
$ws.Range("S13").Value2 = 0.0009984092830520915
$ws.Range("T13").Value2 = 0.0009984092830520915

$ws.Range("G14").Value2 = 13.96350433333333
$ws.Range("H14").Value2 = 41.890513
$ws.Range("I14").Value2 = 0.3275779092170975
$ws.Range("J14").Value2 = 0.3275779092170975
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 1.480335666666667
$ws.Range("N14").Value2 = 4.441007
$ws.Range("O14").Value2 = 0.1826408776454046
$ws.Range("P14").Value2 = 0.1826408776454046
$ws.Range("Q14").Value2 = 20.67067349628789
$ws.Range("R14").Value2 = 186.036061466591
$ws.Range("S14").Value2 = 0.05982911683665735
$ws.Range("T14").Value2 = 0.05982911683665734

$ws.Range("G15").Value2 = 13.96350433333333
$ws.Range("H15").Value2 = 41.890513
$ws.Range("I15").Value2 = 0.3275779092170975
$ws.Range("J15").Value2 = 0.3275779092170975
$ws.Range("K15").Value2 = 3
$ws.Range("L15").Value2 = 1
$ws.Range("M15").Value2 = 0.6291593333333333
$ws.Range("N15").Value2 = 1.887478
$ws.Range("O15").Value2 = 0.07762443032771463
$ws.Range("P15").Value2 = 0.07762443032771463
$ws.Range("Q15").Value2 = 8.78526907735711
$ws.Range("R15").Value2 = 79.067421696214
$ws.Range("S15").Value2 = 0.02542804859092101
$ws.Range("T15").Value2 = 0.02542804859092101

$ws.Range("G16").Value2 = 13.96350433333333
$ws.Range("H16").Value2 = 41.890513
$ws.Range("I16").Value2 = 0.3275779092170975
$ws.Range("J16").Value2 = 0.3275779092170975
$ws.Range("K16").Value2 = 3
$ws.Range("L16").Value2 = 1
$ws.Range("M16").Value2 = 5.953764333333333
$ws.Range("N16").Value2 = 17.861293
$ws.Range("O16").Value2 = 0.734563631492074
$ws.Range("P16").Value2 = 0.734563631492074
$ws.Range("Q16").Value2 = 83.13541406814544
$ws.Range("R16").Value2 = 748.218726613309
$ws.Range("S16").Value2 = 0.2406268185910921
$ws.Range("T16").Value2 = 0.2406268185910921

$ws.Range("G17").Value2 = 13.96350433333333
$ws.Range("H17").Value2 = 41.890513
$ws.Range("I17").Value2 = 0.3275779092170975
$ws.Range("J17").Value2 = 0.3275779092170975
$ws.Range("K17").Value2 = 1
$ws.Range("L17").Value2 = 0.3333333333333333
$ws.Range("M17").Value2 = 0.04191233333333333
$ws.Range("N17").Value2 = 0.125737
$ws.Range("O17").Value2 = 0.005171060534806686
$ws.Range("P17").Value2 = 0.005171060534806686
$ws.Range("Q17").Value2 = 0.5852430481201111
$ws.Range("R17").Value2 = 5.267187433080999
$ws.Range("S17").Value2 = 0.00169392519842702
$ws.Range("T17").Value2 = 0.00169392519842702
